$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 8-10 (old MuSCs -> C3/Itgam rows; data was
# recomputed and consolidated into fewer rows using new TPM-based values)
$ws.Rows("8:10").Delete()

$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "C3"
$ws.Range("C2").Value2 = "Itgam"
$ws.Range("D2").Value2 = "FAPs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.04155
$ws.Range("H2").Value2 = 0.12465
$ws.Range("I2").Value2 = 0.0001466168179836329
$ws.Range("J2").Value2 = 0.0001466168179836329
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.005673666666666667
$ws.Range("N2").Value2 = 0.017021
$ws.Range("O2").Value2 = 0.1234828534325781
$ws.Range("P2").Value2 = 0.1234828534325781
$ws.Range("Q2").Value2 = 0.00023574085
$ws.Range("R2").Value2 = 0.00212166765
$ws.Range("S2").Value2 = [double]"1.810466304582392E-05"
$ws.Range("T2").Value2 = [double]"1.810466304582392E-05"
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "C3"
$ws.Range("C3").Value2 = "Itgam"
$ws.Range("D3").Value2 = "MuSCs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.04155
$ws.Range("H3").Value2 = 0.12465
$ws.Range("I3").Value2 = 0.0001466168179836329
$ws.Range("J3").Value2 = 0.0001466168179836329
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0.3333333333333333
$ws.Range("M3").Value2 = 0.04027333333333333
$ws.Range("N3").Value2 = 0.12082
$ws.Range("O3").Value2 = 0.8765171465674219
$ws.Range("P3").Value2 = 0.876517146567422
$ws.Range("Q3").Value2 = 0.001673357
$ws.Range("R3").Value2 = 0.015060213
$ws.Range("S3").Value2 = 0.0001285121549378089
$ws.Range("T3").Value2 = 0.000128512154937809
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("B4").Value2 = "C3"
$ws.Range("C4").Value2 = "Itgam"
$ws.Range("D4").Value2 = "FAPs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 283.1772156666667
$ws.Range("H4").Value2 = 849.531647
$ws.Range("I4").Value2 = 0.9992428949822291
$ws.Range("J4").Value2 = 0.9992428949822291
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.005673666666666667
$ws.Range("N4").Value2 = 0.017021
$ws.Range("O4").Value2 = 0.1234828534325781
$ws.Range("P4").Value2 = 0.1234828534325781
$ws.Range("Q4").Value2 = 1.606653129287444
$ws.Range("R4").Value2 = 14.459878163587
$ws.Range("S4").Value2 = 0.1233893639446356
$ws.Range("T4").Value2 = 0.1233893639446356
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "C3"
$ws.Range("C5").Value2 = "Itgam"
$ws.Range("D5").Value2 = "MuSCs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 283.1772156666667
$ws.Range("H5").Value2 = 849.531647
$ws.Range("I5").Value2 = 0.9992428949822291
$ws.Range("J5").Value2 = 0.9992428949822291
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.04027333333333333
$ws.Range("N5").Value2 = 0.12082
$ws.Range("O5").Value2 = 0.8765171465674219
$ws.Range("P5").Value2 = 0.876517146567422
$ws.Range("Q5").Value2 = 11.40449039894889
$ws.Range("R5").Value2 = 102.64041359054
$ws.Range("S5").Value2 = 0.8758535310375934
$ws.Range("T5").Value2 = 0.8758535310375936
$ws.Range("A6").Value2 = "MuSCs"
$ws.Range("B6").Value2 = "C3"
$ws.Range("C6").Value2 = "Itgam"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 0.1730073333333333
$ws.Range("H6").Value2 = 0.519022
$ws.Range("I6").Value2 = 0.0006104881997874136
$ws.Range("J6").Value2 = 0.0006104881997874135
$ws.Range("K6").Value2 = 1
$ws.Range("L6").Value2 = 0.3333333333333333
$ws.Range("M6").Value2 = 0.005673666666666667
$ws.Range("N6").Value2 = 0.017021
$ws.Range("O6").Value2 = 0.1234828534325781
$ws.Range("P6").Value2 = 0.1234828534325781
$ws.Range("Q6").Value2 = 0.0009815859402222222
$ws.Range("R6").Value2 = 0.008834273462000001
$ws.Range("S6").Value2 = [double]"7.538482489666765E-05"
$ws.Range("T6").Value2 = [double]"7.538482489666765E-05"
$ws.Range("A7").Value2 = "MuSCs"
$ws.Range("B7").Value2 = "C3"
$ws.Range("C7").Value2 = "Itgam"
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 0.1730073333333333
$ws.Range("H7").Value2 = 0.519022
$ws.Range("I7").Value2 = 0.0006104881997874136
$ws.Range("J7").Value2 = 0.0006104881997874135
$ws.Range("K7").Value2 = 1
$ws.Range("L7").Value2 = 0.3333333333333333
$ws.Range("M7").Value2 = 0.04027333333333333
$ws.Range("N7").Value2 = 0.12082
$ws.Range("O7").Value2 = 0.8765171465674219
$ws.Range("P7").Value2 = 0.876517146567422
$ws.Range("Q7").Value2 = 0.006967582004444444
$ws.Range("R7").Value2 = 0.06270823804
$ws.Range("S7").Value2 = 0.0005351033748907459
$ws.Range("T7").Value2 = 0.0005351033748907458
